# Update the Vietnamese version note: mark "1.1 Background" and
# "1.2 Related work" rows as written by "Thao" instead of "Hao",
# then move the active selection to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Thao"
$ws.Range("C5").Value = "Thao"

$ws.Range("C8").Select()
